# Dispatch Priority by Elec Source.xlsx - apply commit changes
# (updated bcrbq, bgdpbes, bpmccs, crbq, dpbes, pmccs, raf)
#
# Substantive changes from the diff:
#  1. DPbES sheet: "hydro" row (row 6) guaranteed-dispatch flags B6:AE6
#     flip from 1 -> 0 for every forecast year.
#  2. The DPbES worksheet becomes the active tab/sheet (was "About").
#  3. Selection/view state updates:
#       - About sheet view loses tabSelected (no longer the active tab).
#       - DPbES sheet view loses its topLeftCell="D1" and old selection
#         (B12:AE14), replaced with a fresh selection anchored at A6.
#  4. About!A4 loses its (redundant/no-op) cell style so it goes back to
#     the workbook's default formatting.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsDPbES = $wb.Worksheets.Item("DPbES")

# --- 1. Flip the "hydro" dispatch-priority row to 0 for all years ---
$wsDPbES.Range("B6:AE6").Value = 0

# --- 4. About!A4: drop the stray/no-op style back to workbook default ---
$wsAbout.Range("A4").ClearFormats()

# --- 2 & 3. Make DPbES the active sheet, set its selection, and clear
#            the About sheet's prior selection/tab state ---
$wsDPbES.Activate()
$wsDPbES.Range("A6").Select()
